# RefineData.xlsx update
# - Halves the material costs (columns C:G, rows 2-10) on both the
#   "무기" (weapon) and "방어구" (armor) sheets.
# - Switches the active sheet/tab from 무기 to 방어구 and updates the
#   remembered selections on each sheet.

$wb = $excel.ActiveWorkbook
$wsWeapon = $wb.Worksheets.Item(1)   # 무기
$wsArmor  = $wb.Worksheets.Item(2)   # 방어구

# New values for columns C,D,E,F,G (파괴석,돌파석,오레하,명예파편,골드) per row
$weaponData = @{
  2  = @(129,4,2,30,200)
  3  = @(129,4,2,30,200)
  4  = @(129,4,2,30,200)
  5  = @(160,5,2,37,200)
  6  = @(160,5,2,37,200)
  7  = @(160,5,2,37,200)
  8  = @(190,5,3,44,200)
  9  = @(190,6,3,44,200)
  10 = @(190,6,3,44,200)
}

$armorData = @{
  2  = @(78,2,1,21,110)
  3  = @(78,2,1,21,110)
  4  = @(78,2,1,21,110)
  5  = @(96,3,2,25,110)
  6  = @(96,3,2,25,110)
  7  = @(96,3,2,25,110)
  8  = @(114,3,2,30,110)
  9  = @(114,4,2,30,110)
  10 = @(114,4,2,30,110)
}

foreach ($r in $weaponData.Keys) {
    $rowvals = $weaponData[$r]
    for ($i = 0; $i -lt $rowvals.Length; $i++) {
        $wsWeapon.Cells.Item($r, 3 + $i).Value = $rowvals[$i]
    }
}

foreach ($r in $armorData.Keys) {
    $rowvals = $armorData[$r]
    for ($i = 0; $i -lt $rowvals.Length; $i++) {
        $wsArmor.Cells.Item($r, 3 + $i).Value = $rowvals[$i]
    }
}

# Update the remembered selection on each sheet
$wsWeapon.Range("F8").Select()
$wsArmor.Range("H6").Select()

# Make the armor sheet the active tab (this also clears tabSelected on
# the weapon sheet and sets it on the armor sheet)
$wsArmor.Activate()
